$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 659.4167
$ws.Range("I43").Value = 753.6667
$ws.Range("J43").Value = 628
$ws.Range("K43").Value = 753.6667
$ws.Range("L43").Value = 628
$ws.Range("M43").Value = -684.6667
$ws.Range("N43").Value = -766
$ws.Range("H51").Value = 1716.5
$ws.Range("I51").Value = 1999.5
$ws.Range("K51").Value = 1999.5
$ws.Range("M51").Value = -1515.5
$ws.Range("H64").Value = 3978.392
$ws.Range("I64").Value = 3833.2917
$ws.Range("J64").Value = 4107.3706
$ws.Range("K64").Value = 3833.2917
$ws.Range("L64").Value = 4107.3706
$ws.Range("M64").Value = -3585.2917
$ws.Range("N64").Value = -4603.3706
$ws.Range("H67").Value = 3978.392
$ws.Range("I67").Value = 3833.2917
$ws.Range("J67").Value = 4107.3706
$ws.Range("K67").Value = 3833.2917
$ws.Range("L67").Value = 4107.3706
$ws.Range("M67").Value = -2975.2917
$ws.Range("N67").Value = -5823.3706
$ws.Range("H88").Value = 4308.636
$ws.Range("I88").Value = 739.1
$ws.Range("K88").Value = 739.1
$ws.Range("M88").Value = -333.1
$ws.Range("H91").Value = 4308.636
$ws.Range("I91").Value = 739.1
$ws.Range("K91").Value = 739.1
$ws.Range("M91").Value = 664.9
$ws.Range("H98").Value = 1993.9474
$ws.Range("I98").Value = 1642
$ws.Range("J98").Value = 2597.2856
$ws.Range("K98").Value = 1642
$ws.Range("L98").Value = 2597.2856
$ws.Range("M98").Value = -144
$ws.Range("N98").Value = -5593.2856
$ws.Range("H116").Value = 2139762.8
$ws.Range("I116").Value = 5497454.5
$ws.Range("J116").Value = 3049.818
$ws.Range("K116").Value = 5497454.5
$ws.Range("L116").Value = 3049.818
$ws.Range("M116").Value = -5494012.5
$ws.Range("N116").Value = -9933.817999999999
$ws.Range("H122").Value = 1993.9474
$ws.Range("I122").Value = 1642
$ws.Range("J122").Value = 2597.2856
$ws.Range("K122").Value = 4926
$ws.Range("L122").Value = 7791.8568
$ws.Range("M122").Value = -2476
$ws.Range("N122").Value = -12691.8568
$ws.Range("H125").Value = 920.6667
$ws.Range("I125").Value = 910.75
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 8196.75
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -5736.75
$ws.Range("N125").Value = -13920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 5270.4443
$ws.Range("I26").Value = 2337.8667
$ws.Range("J26").Value = 19933.334
$ws.Range("K26").Value = 2337.8667
$ws.Range("L26").Value = 19933.334
$ws.Range("M26").Value = -2007.8667
$ws.Range("N26").Value = -20593.334
$ws.Range("H32").Value = 20004564
$ws.Range("I32").Value = 25002380
$ws.Range("J32").Value = 13303.1
$ws.Range("K32").Value = 25002380
$ws.Range("L32").Value = 13303.1
$ws.Range("M32").Value = -25002093
$ws.Range("N32").Value = -13877.1
$ws.Range("H122").Value = 1882.4
$ws.Range("I122").Value = 1603
$ws.Range("K122").Value = 4809
$ws.Range("M122").Value = -2359
$ws.Range("H130").Value = 29142.666
$ws.Range("J130").Value = 29142.666
$ws.Range("L130").Value = 29142.666
$ws.Range("N130").Value = -39182.666
$ws.Range("H133").Value = 36490
$ws.Range("J133").Value = 36490
$ws.Range("L133").Value = 36490
$ws.Range("N133").Value = -41550
$ws.Range("H135").Value = 63479.57
$ws.Range("J135").Value = 63479.57
$ws.Range("L135").Value = 63479.57
$ws.Range("N135").Value = -73619.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1100.6666
$ws.Range("I22").Value = 900.5
$ws.Range("J22").Value = 1501
$ws.Range("K22").Value = 900.5
$ws.Range("L22").Value = 1501
$ws.Range("M22").Value = -550.5
$ws.Range("N22").Value = -2201
$ws.Range("H94").Value = 4465.1763
$ws.Range("I94").Value = 2603.6667
$ws.Range("J94").Value = 4864.0713
$ws.Range("K94").Value = 2603.6667
$ws.Range("L94").Value = 4864.0713
$ws.Range("M94").Value = -2152.6667
$ws.Range("N94").Value = -5766.0713
$ws.Range("H130").Value = 43695
$ws.Range("J130").Value = 43695
$ws.Range("L130").Value = 43695
$ws.Range("N130").Value = -53735
$ws.Range("H132").Value = 2575.5881
$ws.Range("I132").Value = 1545.1538
$ws.Range("K132").Value = 4635.4614
$ws.Range("M132").Value = -2105.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 222.41176
$ws.Range("I23").Value = 196.75
$ws.Range("K23").Value = 590.25
$ws.Range("M23").Value = -355.25
$ws.Range("H97").Value = 553.9
$ws.Range("I97").Value = 455.83334
$ws.Range("J97").Value = 701
$ws.Range("K97").Value = 1367.50002
$ws.Range("L97").Value = 2103
$ws.Range("M97").Value = -871.5000199999999
$ws.Range("N97").Value = -3095
$ws.Range("H98").Value = 251
$ws.Range("I98").Value = 251
$ws.Range("K98").Value = 753
$ws.Range("M98").Value = 745
$ws.Range("H122").Value = 8925.73
$ws.Range("I122").Value = 12220.556
$ws.Range("K122").Value = 109985.004
$ws.Range("M122").Value = -107535.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 23340
$ws.Range("I113").Value = 2900
$ws.Range("K113").Value = 2900
$ws.Range("M113").Value = -730
$ws.Range("H127").Value = 29000
$ws.Range("J127").Value = 29000
$ws.Range("L127").Value = 29000
$ws.Range("N127").Value = -38920

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4742.3887
$ws.Range("I40").Value = 4247.357
$ws.Range("J40").Value = 6475
$ws.Range("K40").Value = 4247.357
$ws.Range("L40").Value = 6475
$ws.Range("M40").Value = -4111.357
$ws.Range("N40").Value = -6747
$ws.Range("H122").Value = 2153.3684
$ws.Range("I122").Value = 2153.3684
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6460.1052
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -4010.1052
$ws.Range("H128").Value = 28229
$ws.Range("J128").Value = 28229
$ws.Range("L128").Value = 28229
$ws.Range("N128").Value = -38189
$ws.Range("H132").Value = 3480.3809
$ws.Range("I132").Value = 3214.0293
$ws.Range("J132").Value = 4612.375
$ws.Range("K132").Value = 9642.0879
$ws.Range("L132").Value = 13837.125
$ws.Range("M132").Value = -7112.0879
$ws.Range("N132").Value = -18897.125
$ws.Range("H138").Value = 22813.666
$ws.Range("J138").Value = 22813.666
$ws.Range("L138").Value = 22813.666
$ws.Range("N138").Value = -33093.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 823.6923
$ws.Range("I100").Value = 438.125
$ws.Range("J100").Value = 995.05554
$ws.Range("K100").Value = 876.25
$ws.Range("L100").Value = 1990.11108
$ws.Range("M100").Value = -335.25
$ws.Range("N100").Value = -3072.11108
$ws.Range("H126").Value = 1910.9286
$ws.Range("I126").Value = 2036.4
$ws.Range("J126").Value = 1597.25
$ws.Range("K126").Value = 6109.200000000001
$ws.Range("L126").Value = 4791.75
$ws.Range("M126").Value = -3639.200000000001
$ws.Range("N126").Value = -9731.75
$ws.Range("H132").Value = 2131.3333
$ws.Range("I132").Value = 1601.2916
$ws.Range("J132").Value = 3544.7778
$ws.Range("K132").Value = 4803.8748
$ws.Range("L132").Value = 10634.3334
$ws.Range("M132").Value = -2273.8748
$ws.Range("N132").Value = -15694.3334
